# Make the doc_institution demo data more realistic and complete:
#  - row 2 (pdf_wiki) gets a real institution_id value ("suisse") instead
#    of the placeholder "institution_1"
#  - a new row 3 is appended (pop_com_1 / ofs-pop)
#  - the Tableau1 table (and its autofilter) grows to include the new row
#  - the active selection moves to the first empty row below the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder institution_id for the existing data row.
$ws.Range("B2").Value = "suisse"

# Append the new data row.
$ws.Range("A3").Value = "pop_com_1"
$ws.Range("B3").Value = "ofs-pop"

# Grow the table to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B3")) | Out-Null

# Leave the selection on the next free row, as in the saved file.
$ws.Range("B4").Select() | Out-Null
